$d = $word.ActiveDocument

# Locate the end of the sentence "...uses the Librarian to access robots and
# their information." in the Description / Overview paragraph.
$findRange = $d.Content
$null = $findRange.Find.Execute("access robots and their information.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Collapse the found range to its end (still inside the paragraph, right before
# the paragraph mark) and append the new sentence there as a new run.
$findRange.Collapse(0)

$newSentence = " There is a selection box to choose a robot, a textbox to view the robots stats, a box to sort the robots in the selection box, and 3 buttons to upload the statistics, reset the statistics, and go back to the main menu."

# Insert the new sentence plus a temporary one-character placeholder. The
# placeholder lets us build a non-degenerate bookmark range and avoid placing a
# zero-width range exactly on the paragraph-end boundary.
$findRange.InsertAfter($newSentence + "X")

# Re-anchor the "_GoBack" bookmark (which already exists elsewhere in the
# document) around the placeholder character; because bookmark names are
# unique, this moves the existing bookmark to this new location.
$placeholderRange = $d.Range($findRange.End - 1, $findRange.End)
$d.Bookmarks.Add("_GoBack", $placeholderRange)

# Delete the placeholder character, leaving the (now collapsed) bookmark
# sitting immediately after the new run and before the paragraph mark.
$placeholderRange = $d.Range($findRange.End - 1, $findRange.End)
$placeholderRange.Delete()
